$wb = $excel.ActiveWorkbook

# Reorder tabs: move review_info before hotel_info
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($hotelInfo)

# Re-fetch a fresh reference to hotel_info now that tab order changed
$hotelInfo = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info between Hotel_Name (B) and City (C)
$hotelInfo.Columns("C:C").Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"
